$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Pin Header JLCPCB part number (row 33)
$ws.Range("D33").Value = "C2905434"

# Update the Audio Jack connector row (row 28)
$ws.Range("C28").Value = "PJ-3200"
$ws.Range("D28").Value = "C2689690"
$ws.Range("A28").Value = "Audio Jack"

# Update the view/selection state to match the saved workbook view
$win = $excel.ActiveWindow
$win.ScrollRow = 26
$win.ScrollColumn = 1
$ws.Range("A28").Select() | Out-Null
